$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 5
$ws.Range("L2").Value = "stimuli/img_inqod.png"
$ws.Range("M2").Value = 70.84848484848484
$ws.Range("N2").Value = 50.63636363636363
$ws.Range("O2").Value = 60.74242424242424
$ws.Range("P2").Value = 33
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 5
$ws.Range("S2").Value = 5
$ws.Range("T2").Value = 5
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 5

# Row 3
$ws.Range("C3").Value = 5
$ws.Range("I3").Value = "target"
$ws.Range("J3").Value = "old"
$ws.Range("K3").Value = "j"
$ws.Range("L3").Value = "stimuli/img_aplao.png"
$ws.Range("M3").Value = 64.0909090909091
$ws.Range("N3").Value = 40.75757575757576
$ws.Range("O3").Value = 52.42424242424242
$ws.Range("P3").Value = 33
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 3
$ws.Range("U3").Value = 3
$ws.Range("V3").Value = 3

# Row 4
$ws.Range("C4").Value = 5
$ws.Range("I4").ClearContents()
$ws.Range("J4").Value = "new"
$ws.Range("K4").Value = "f"
$ws.Range("L4").Value = "stimuli/img_mgnmm.png"
$ws.Range("M4").Value = 79.1470588235294
$ws.Range("N4").Value = 60.38235294117647
$ws.Range("O4").Value = 69.76470588235294
$ws.Range("P4").Value = 34
$ws.Range("Q4").Value = 8
$ws.Range("R4").Value = 8
$ws.Range("S4").Value = 8
$ws.Range("T4").Value = 8
$ws.Range("U4").Value = 8
$ws.Range("V4").Value = 8

# Row 5
$ws.Range("C5").Value = 5
$ws.Range("L5").Value = "stimuli/img_30vz5.png"
$ws.Range("M5").Value = 86.21212121212122
$ws.Range("N5").Value = 68.27272727272727
$ws.Range("O5").Value = 77.24242424242425
$ws.Range("P5").Value = 33
$ws.Range("Q5").Value = 10
$ws.Range("R5").Value = 10
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = 10
$ws.Range("U5").Value = 10
$ws.Range("V5").Value = 10

# Row 6
$ws.Range("C6").Value = 5
$ws.Range("L6").Value = "stimuli/img_esb4r.png"
$ws.Range("M6").Value = 60.73529411764706
$ws.Range("N6").Value = 38.58823529411764
$ws.Range("O6").Value = 49.66176470588235
$ws.Range("P6").Value = 34
$ws.Range("Q6").Value = 3
$ws.Range("R6").Value = 3
$ws.Range("S6").Value = 3
$ws.Range("T6").Value = 3
$ws.Range("U6").Value = 3
$ws.Range("V6").Value = 3

# Row 7
$ws.Range("C7").Value = 5
$ws.Range("I7").ClearContents()
$ws.Range("J7").Value = "new"
$ws.Range("K7").Value = "f"
$ws.Range("L7").Value = "stimuli/img_mjxmq.png"
$ws.Range("M7").Value = 77.07692307692308
$ws.Range("N7").Value = 58.15384615384615
$ws.Range("O7").Value = 67.61538461538461
$ws.Range("P7").Value = 39
$ws.Range("Q7").Value = 7
$ws.Range("R7").Value = 7
$ws.Range("S7").Value = 7
$ws.Range("T7").Value = 7
$ws.Range("U7").Value = 7
$ws.Range("V7").Value = 7

# Row 8
$ws.Range("C8").Value = 5
$ws.Range("L8").Value = "stimuli/img_qmgwq.png"
$ws.Range("M8").Value = 84.58333333333333
$ws.Range("N8").Value = 64.44444444444444
$ws.Range("O8").Value = 74.51388888888889
$ws.Range("P8").Value = 36
$ws.Range("Q8").Value = 9
$ws.Range("R8").Value = 9
$ws.Range("S8").Value = 9
$ws.Range("T8").Value = 9
$ws.Range("U8").Value = 9
$ws.Range("V8").Value = 9

# Row 9
$ws.Range("C9").Value = 5
$ws.Range("I9").ClearContents()
$ws.Range("J9").Value = "new"
$ws.Range("K9").Value = "f"
$ws.Range("L9").Value = "stimuli/img_z293c.png"
$ws.Range("M9").Value = 71.26470588235294
$ws.Range("N9").Value = 46.88235294117647
$ws.Range("O9").Value = 59.07352941176471
$ws.Range("P9").Value = 34
$ws.Range("Q9").Value = 5
$ws.Range("R9").Value = 5
$ws.Range("S9").Value = 5
$ws.Range("T9").Value = 5
$ws.Range("U9").Value = 5
$ws.Range("V9").Value = 5

# Row 10
$ws.Range("C10").Value = 5
$ws.Range("L10").Value = "stimuli/img_vbrb7.png"
$ws.Range("M10").Value = 85.5625
$ws.Range("N10").Value = 71.46875
$ws.Range("O10").Value = 78.515625
$ws.Range("P10").Value = 32
$ws.Range("Q10").Value = 10
$ws.Range("R10").Value = 10
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = 10
$ws.Range("U10").Value = 10
$ws.Range("V10").Value = 10

# Row 11
$ws.Range("C11").Value = 5
$ws.Range("I11").Value = "target"
$ws.Range("J11").Value = "old"
$ws.Range("K11").Value = "j"
$ws.Range("L11").Value = "stimuli/img_iyxnj.png"
$ws.Range("M11").Value = 75.30555555555556
$ws.Range("N11").Value = 54.33333333333334
$ws.Range("O11").Value = 64.81944444444444
$ws.Range("P11").Value = 36
$ws.Range("Q11").Value = 6
$ws.Range("R11").Value = 6
$ws.Range("S11").Value = 6
$ws.Range("T11").Value = 6
$ws.Range("U11").Value = 6
$ws.Range("V11").Value = 6

# Row 12
$ws.Range("C12").Value = 5
$ws.Range("L12").Value = "stimuli/img_ce9vx.png"
$ws.Range("M12").Value = 75.9090909090909
$ws.Range("N12").Value = 57.12121212121212
$ws.Range("O12").Value = 66.51515151515152
$ws.Range("P12").Value = 33
$ws.Range("Q12").Value = 7
$ws.Range("R12").Value = 7
$ws.Range("S12").Value = 7
$ws.Range("T12").Value = 7
$ws.Range("U12").Value = 7
$ws.Range("V12").Value = 7

# Row 13
$ws.Range("C13").Value = 5
$ws.Range("L13").Value = "stimuli/img_cv6mf.png"
$ws.Range("M13").Value = 66.8
$ws.Range("N13").Value = 42.08
$ws.Range("O13").Value = 54.44
$ws.Range("P13").Value = 25
$ws.Range("Q13").Value = 4
$ws.Range("R13").Value = 4
$ws.Range("S13").Value = 4
$ws.Range("T13").Value = 4
$ws.Range("U13").Value = 4
$ws.Range("V13").Value = 4

# Row 14
$ws.Range("C14").Value = 5
$ws.Range("I14").Value = "target"
$ws.Range("J14").Value = "old"
$ws.Range("K14").Value = "j"
$ws.Range("L14").Value = "stimuli/img_p3hpc.png"
$ws.Range("M14").Value = 72.83333333333333
$ws.Range("N14").Value = 52.22222222222222
$ws.Range("O14").Value = 62.52777777777777
$ws.Range("P14").Value = 36
$ws.Range("Q14").Value = 6
$ws.Range("R14").Value = 6
$ws.Range("S14").Value = 6
$ws.Range("T14").Value = 6
$ws.Range("U14").Value = 6
$ws.Range("V14").Value = 6

# Row 15
$ws.Range("C15").Value = 5
$ws.Range("L15").Value = "stimuli/img_lszzj.png"
$ws.Range("M15").Value = 64.70588235294117
$ws.Range("N15").Value = 45.58823529411764
$ws.Range("O15").Value = 55.14705882352941
$ws.Range("P15").Value = 34
$ws.Range("Q15").Value = 4
$ws.Range("R15").Value = 4
$ws.Range("S15").Value = 4
$ws.Range("T15").Value = 4
$ws.Range("U15").Value = 4
$ws.Range("V15").Value = 4

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("I16").ClearContents()
$ws.Range("J16").Value = "new"
$ws.Range("K16").Value = "f"
$ws.Range("L16").Value = "stimuli/img_7ed9m.png"
$ws.Range("M16").Value = 80.71875
$ws.Range("N16").Value = 58.65625
$ws.Range("O16").Value = 69.6875
$ws.Range("P16").Value = 32
$ws.Range("Q16").Value = 8
$ws.Range("R16").Value = 8
$ws.Range("S16").Value = 8
$ws.Range("T16").Value = 8
$ws.Range("U16").Value = 8
$ws.Range("V16").Value = 8

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("I17").ClearContents()
$ws.Range("J17").Value = "new"
$ws.Range("K17").Value = "f"
$ws.Range("L17").Value = "stimuli/img_r2lxk.png"
$ws.Range("M17").Value = 89.24242424242425
$ws.Range("N17").Value = 67.6969696969697
$ws.Range("O17").Value = 78.46969696969697
$ws.Range("P17").Value = 33
$ws.Range("Q17").Value = 10
$ws.Range("R17").Value = 10
$ws.Range("S17").Value = 10
$ws.Range("T17").Value = 10
$ws.Range("U17").Value = 10
$ws.Range("V17").Value = 10

# Row 18
$ws.Range("C18").Value = 5
$ws.Range("L18").Value = "stimuli/img_p659z.png"
$ws.Range("M18").Value = 84.21621621621621
$ws.Range("N18").Value = 65.37837837837837
$ws.Range("O18").Value = 74.79729729729729
$ws.Range("P18").Value = 37
$ws.Range("Q18").Value = 9
$ws.Range("R18").Value = 9
$ws.Range("S18").Value = 9
$ws.Range("T18").Value = 9
$ws.Range("U18").Value = 9
$ws.Range("V18").Value = 9

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("L19").Value = "stimuli/img_ewrjk.png"
$ws.Range("M19").Value = 73.0909090909091
$ws.Range("N19").Value = 53.39393939393939
$ws.Range("O19").Value = 63.24242424242424
$ws.Range("P19").Value = 33
$ws.Range("Q19").Value = 6
$ws.Range("R19").Value = 6
$ws.Range("S19").Value = 6
$ws.Range("T19").Value = 6
$ws.Range("U19").Value = 6
$ws.Range("V19").Value = 6

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("L20").Value = "stimuli/img_t90e2.png"
$ws.Range("M20").Value = 83.0625
$ws.Range("N20").Value = 61.96875
$ws.Range("O20").Value = 72.515625
$ws.Range("P20").Value = 32
$ws.Range("Q20").Value = 9
$ws.Range("R20").Value = 9
$ws.Range("S20").Value = 9
$ws.Range("T20").Value = 9
$ws.Range("U20").Value = 9
$ws.Range("V20").Value = 9

# Row 21
$ws.Range("C21").Value = 5
$ws.Range("H21").ClearContents()
$ws.Range("J21").Value = "catch"
$ws.Range("K21").Value = "f"
$ws.Range("L21").Value = "stimuli/catch_28.jpg"
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("O21").ClearContents()
$ws.Range("P21").ClearContents()
$ws.Range("Q21").ClearContents()
$ws.Range("R21").ClearContents()
$ws.Range("S21").ClearContents()
$ws.Range("T21").ClearContents()
$ws.Range("U21").ClearContents()
$ws.Range("V21").ClearContents()

# Row 22
$ws.Range("C22").Value = 5
$ws.Range("I22").ClearContents()
$ws.Range("J22").Value = "new"
$ws.Range("K22").Value = "f"
$ws.Range("L22").Value = "stimuli/img_kwxq1.png"
$ws.Range("M22").Value = 68.53125
$ws.Range("N22").Value = 44.09375
$ws.Range("O22").Value = 56.3125
$ws.Range("P22").Value = 32
$ws.Range("Q22").Value = 4
$ws.Range("R22").Value = 4
$ws.Range("S22").Value = 4
$ws.Range("T22").Value = 4
$ws.Range("U22").Value = 4
$ws.Range("V22").Value = 4

# Row 23
$ws.Range("C23").Value = 5
$ws.Range("I23").ClearContents()
$ws.Range("J23").Value = "new"
$ws.Range("K23").Value = "f"
$ws.Range("L23").Value = "stimuli/img_wppku.png"
$ws.Range("M23").Value = 75.02941176470588
$ws.Range("N23").Value = 53.05882352941177
$ws.Range("O23").Value = 64.04411764705883
$ws.Range("P23").Value = 34
$ws.Range("Q23").Value = 6
$ws.Range("R23").Value = 6
$ws.Range("S23").Value = 6
$ws.Range("T23").Value = 6
$ws.Range("U23").Value = 6
$ws.Range("V23").Value = 6

# Row 24
$ws.Range("C24").Value = 5
$ws.Range("L24").Value = "stimuli/img_7ucnr.png"
$ws.Range("M24").Value = 70.39393939393939
$ws.Range("N24").Value = 47.90909090909091
$ws.Range("O24").Value = 59.15151515151515
$ws.Range("P24").Value = 33
$ws.Range("Q24").Value = 5
$ws.Range("R24").Value = 5
$ws.Range("S24").Value = 5
$ws.Range("T24").Value = 5
$ws.Range("U24").Value = 5
$ws.Range("V24").Value = 5

# Row 25
$ws.Range("C25").Value = 5
$ws.Range("I25").Value = "target"
$ws.Range("J25").Value = "old"
$ws.Range("K25").Value = "j"
$ws.Range("L25").Value = "stimuli/img_7wul8.png"
$ws.Range("M25").Value = 43.03030303030303
$ws.Range("N25").Value = 25.54545454545455
$ws.Range("O25").Value = 34.28787878787879
$ws.Range("P25").Value = 33
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = 1
$ws.Range("S25").Value = 1
$ws.Range("T25").Value = 1
$ws.Range("U25").Value = 1
$ws.Range("V25").Value = 1

# Row 26
$ws.Range("C26").Value = 5
$ws.Range("L26").Value = "stimuli/img_7w5tw.png"
$ws.Range("M26").Value = 53.2258064516129
$ws.Range("N26").Value = 28.90322580645161
$ws.Range("O26").Value = 41.06451612903226
$ws.Range("P26").Value = 31
$ws.Range("Q26").Value = 2
$ws.Range("R26").Value = 2
$ws.Range("S26").Value = 2
$ws.Range("T26").Value = 2
$ws.Range("U26").Value = 2
$ws.Range("V26").Value = 2

# Row 27
$ws.Range("C27").Value = 5
$ws.Range("H27").Value = "kitchens"
$ws.Range("J27").Value = "new"
$ws.Range("K27").Value = "f"
$ws.Range("L27").Value = "stimuli/img_1ao2d.png"
$ws.Range("M27").Value = 38.77777777777778
$ws.Range("N27").Value = 18.75
$ws.Range("O27").Value = 28.76388888888889
$ws.Range("P27").Value = 36
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = 1
$ws.Range("S27").Value = 1
$ws.Range("T27").Value = 1
$ws.Range("U27").Value = 1
$ws.Range("V27").Value = 1

# Row 28
$ws.Range("C28").Value = 5
$ws.Range("L28").Value = "stimuli/img_yeh72.png"
$ws.Range("M28").Value = 68.66666666666667
$ws.Range("N28").Value = 45.21212121212121
$ws.Range("O28").Value = 56.93939393939394
$ws.Range("P28").Value = 33
$ws.Range("Q28").Value = 4
$ws.Range("R28").Value = 4
$ws.Range("S28").Value = 4
$ws.Range("T28").Value = 4
$ws.Range("U28").Value = 4
$ws.Range("V28").Value = 4

# Row 29
$ws.Range("C29").Value = 5
$ws.Range("L29").Value = "stimuli/img_eatdk.png"
$ws.Range("M29").Value = 81.40625
$ws.Range("N29").Value = 61.375
$ws.Range("O29").Value = 71.390625
$ws.Range("P29").Value = 32
$ws.Range("Q29").Value = 8
$ws.Range("R29").Value = 8
$ws.Range("S29").Value = 8
$ws.Range("T29").Value = 8
$ws.Range("U29").Value = 8
$ws.Range("V29").Value = 8

# Row 30
$ws.Range("C30").Value = 5
$ws.Range("I30").ClearContents()
$ws.Range("J30").Value = "old"
$ws.Range("K30").Value = "j"
$ws.Range("L30").Value = "stimuli/img_wyl6z.png"
$ws.Range("M30").Value = 59.8235294117647
$ws.Range("N30").Value = 36.23529411764706
$ws.Range("O30").Value = 48.02941176470588
$ws.Range("P30").Value = 34
$ws.Range("Q30").Value = 3
$ws.Range("R30").Value = 3
$ws.Range("S30").Value = 3
$ws.Range("T30").Value = 3
$ws.Range("U30").Value = 3
$ws.Range("V30").Value = 3

# Row 31
$ws.Range("C31").Value = 5
$ws.Range("I31").Value = "target"
$ws.Range("J31").Value = "old"
$ws.Range("K31").Value = "j"
$ws.Range("L31").Value = "stimuli/img_nyv2b.png"
$ws.Range("M31").Value = 11.91176470588235
$ws.Range("N31").Value = 6.852941176470588
$ws.Range("O31").Value = 9.382352941176471
$ws.Range("P31").Value = 34
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = 1
$ws.Range("S31").Value = 1
$ws.Range("T31").Value = 1
$ws.Range("U31").Value = 1
$ws.Range("V31").Value = 1

# Row 32
$ws.Range("C32").Value = 5
$ws.Range("L32").Value = "stimuli/img_mawe6.png"
$ws.Range("M32").Value = 83.48387096774194
$ws.Range("N32").Value = 65.54838709677419
$ws.Range("O32").Value = 74.51612903225806
$ws.Range("P32").Value = 31
$ws.Range("Q32").Value = 9
$ws.Range("R32").Value = 9
$ws.Range("S32").Value = 9
$ws.Range("T32").Value = 9
$ws.Range("U32").Value = 9
$ws.Range("V32").Value = 9

# Row 33
$ws.Range("C33").Value = 5
$ws.Range("L33").Value = "stimuli/img_njmgp.png"
$ws.Range("M33").Value = 80.48148148148148
$ws.Range("N33").Value = 58.4074074074074
$ws.Range("O33").Value = 69.44444444444444
$ws.Range("P33").Value = 27
$ws.Range("Q33").Value = 8
$ws.Range("R33").Value = 8
$ws.Range("S33").Value = 8
$ws.Range("T33").Value = 8
$ws.Range("U33").Value = 8
$ws.Range("V33").Value = 8

# Row 34
$ws.Range("C34").Value = 5
$ws.Range("I34").Value = "target"
$ws.Range("J34").Value = "old"
$ws.Range("K34").Value = "j"
$ws.Range("L34").Value = "stimuli/img_cnyac.png"
$ws.Range("M34").Value = 69.1470588235294
$ws.Range("N34").Value = 47.8235294117647
$ws.Range("O34").Value = 58.48529411764706
$ws.Range("P34").Value = 34
$ws.Range("Q34").Value = 5
$ws.Range("R34").Value = 5
$ws.Range("S34").Value = 5
$ws.Range("T34").Value = 5
$ws.Range("U34").Value = 5
$ws.Range("V34").Value = 5

# Row 35
$ws.Range("C35").Value = 5
$ws.Range("I35").ClearContents()
$ws.Range("J35").Value = "new"
$ws.Range("K35").Value = "f"
$ws.Range("L35").Value = "stimuli/img_ikk62.png"
$ws.Range("M35").Value = 37.48780487804878
$ws.Range("N35").Value = 21.07317073170732
$ws.Range("O35").Value = 29.28048780487805
$ws.Range("P35").Value = 41
$ws.Range("Q35").Value = 1
$ws.Range("R35").Value = 1
$ws.Range("S35").Value = 1
$ws.Range("T35").Value = 1
$ws.Range("U35").Value = 1
$ws.Range("V35").Value = 1

# Row 36
$ws.Range("C36").Value = 5
$ws.Range("L36").Value = "stimuli/img_89rmb.png"
$ws.Range("M36").Value = 55.18518518518518
$ws.Range("N36").Value = 29.25925925925926
$ws.Range("O36").Value = 42.22222222222222
$ws.Range("P36").Value = 27
$ws.Range("Q36").Value = 2
$ws.Range("R36").Value = 2
$ws.Range("S36").Value = 2
$ws.Range("T36").Value = 2
$ws.Range("U36").Value = 2
$ws.Range("V36").Value = 2

# Row 37
$ws.Range("C37").Value = 5
$ws.Range("I37").Value = "target"
$ws.Range("J37").Value = "old"
$ws.Range("K37").Value = "j"
$ws.Range("L37").Value = "stimuli/img_d8xbu.png"
$ws.Range("M37").Value = 91.36363636363636
$ws.Range("N37").Value = 73.18181818181819
$ws.Range("O37").Value = 82.27272727272728
$ws.Range("P37").Value = 33
$ws.Range("Q37").Value = 10
$ws.Range("R37").Value = 10
$ws.Range("S37").Value = 10
$ws.Range("T37").Value = 10
$ws.Range("U37").Value = 10
$ws.Range("V37").Value = 10

# Row 38
$ws.Range("C38").Value = 5
$ws.Range("I38").ClearContents()
$ws.Range("J38").Value = "new"
$ws.Range("K38").Value = "f"
$ws.Range("L38").Value = "stimuli/img_xdhz2.png"
$ws.Range("M38").Value = 63.3
$ws.Range("N38").Value = 37.25
$ws.Range("O38").Value = 50.275
$ws.Range("P38").Value = 40
$ws.Range("Q38").Value = 3
$ws.Range("R38").Value = 3
$ws.Range("S38").Value = 3
$ws.Range("T38").Value = 3
$ws.Range("U38").Value = 3
$ws.Range("V38").Value = 3

# Row 39
$ws.Range("C39").Value = 5
$ws.Range("I39").Value = "target"
$ws.Range("J39").Value = "old"
$ws.Range("K39").Value = "j"
$ws.Range("L39").Value = "stimuli/img_ye5sl.png"
$ws.Range("M39").Value = 53.2258064516129
$ws.Range("N39").Value = 34.45161290322581
$ws.Range("O39").Value = 43.83870967741936
$ws.Range("P39").Value = 31
$ws.Range("Q39").Value = 2
$ws.Range("R39").Value = 2
$ws.Range("S39").Value = 2
$ws.Range("T39").Value = 2
$ws.Range("U39").Value = 2
$ws.Range("V39").Value = 2

# Row 40
$ws.Range("C40").Value = 5
$ws.Range("L40").Value = "stimuli/img_es7o2.png"
$ws.Range("M40").Value = 52.48571428571429
$ws.Range("N40").Value = 27.54285714285714
$ws.Range("O40").Value = 40.01428571428572
$ws.Range("P40").Value = 35
$ws.Range("Q40").Value = 2
$ws.Range("R40").Value = 2
$ws.Range("S40").Value = 2
$ws.Range("T40").Value = 2
$ws.Range("U40").Value = 2
$ws.Range("V40").Value = 2

# Row 41
$ws.Range("C41").Value = 5
$ws.Range("I41").ClearContents()
$ws.Range("J41").Value = "new"
$ws.Range("K41").Value = "f"
$ws.Range("L41").Value = "stimuli/img_zi8qc.png"
$ws.Range("M41").Value = 77.14285714285714
$ws.Range("N41").Value = 57.02857142857143
$ws.Range("O41").Value = 67.08571428571429
$ws.Range("P41").Value = 35
$ws.Range("Q41").Value = 7
$ws.Range("R41").Value = 7
$ws.Range("S41").Value = 7
$ws.Range("T41").Value = 7
$ws.Range("U41").Value = 7
$ws.Range("V41").Value = 7

# Row 42
$ws.Range("C42").Value = 5
$ws.Range("I42").Value = "target"
$ws.Range("J42").Value = "old"
$ws.Range("K42").Value = "j"
$ws.Range("L42").Value = "stimuli/img_6nbgt.png"
$ws.Range("M42").Value = 78.45161290322581
$ws.Range("N42").Value = 57.83870967741935
$ws.Range("O42").Value = 68.14516129032258
$ws.Range("P42").Value = 31
$ws.Range("Q42").Value = 7
$ws.Range("R42").Value = 7
$ws.Range("S42").Value = 7
$ws.Range("T42").Value = 7
$ws.Range("U42").Value = 7
$ws.Range("V42").Value = 7
